$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new row of data (row 22): Leetcode Question No. / Question
$ws.Range("A22").Value = "328/GFG"
$ws.Range("B22").Value = "Odd Even Linked List"

# Move the selection to reflect the post-edit cursor position (B23)
$ws.Range("B23").Select()
